$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "...outputs are described in step 13" -> split into two runs, and
#    change the step number from 13 to 11:
#      run A: "...and the "
#      run B: "outputs are described in step 11"
# ---------------------------------------------------------------------

# First fix the digits while the sentence is still a single run, so the
# run-split we do next lands on the final ("step 11") text.
$d.Content.Find.Execute( `
    "outputs are described in step 13", $false, $false, $false, $false, `
    $false, $true, 1, $false, "outputs are described in step 11", 2) | Out-Null

# Now re-find that (corrected) phrase and force a run break right before
# it by nudging a character-formatting property on just that sub-range.
$splitRng = $d.Content
$splitRng.Find.Execute( `
    "outputs are described in step 11", $false, $false, $false, $false, `
    $false, $true, 1, $false, "", 0) | Out-Null
$splitRng.Font.Bold = 1
$splitRng.Font.Bold = 0

# ---------------------------------------------------------------------
# 2) Remove the "## Source of data" heading paragraph and relocate the
#    automatic "_GoBack" bookmark from its old spot (after the
#    "pre-process text ..." paragraph) to the start of the paragraph
#    that now immediately follows the deleted heading ("The dataset was
#    obtained ...").
# ---------------------------------------------------------------------

$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

$headingRng = $d.Content
$headingRng.Find.Execute( `
    "## Source of data", $false, $false, $false, $false, $false, $true, `
    1, $false, "", 0) | Out-Null
$headingPara = $headingRng.Paragraphs(1)
$headingPara.Range.Delete()

$targetRng = $d.Content
$targetRng.Find.Execute( `
    "The dataset was obtained", $false, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null
$targetRng.Collapse(1)
$d.Bookmarks.Add("_GoBack", $targetRng)

# ---------------------------------------------------------------------
# 3) Update the cached page-number field result in the default footer
#    from "8" to "1".
# ---------------------------------------------------------------------

$footer = $d.Sections.Item(1).Footers.Item(1)
$footer.Range.Find.Execute( `
    "8", $false, $false, $false, $false, $false, $true, 1, $false, "1", 2) | Out-Null
